$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ticket IDs (column A) and status (column B) for the active rows.
# New tickets shift into rows 1-10 (the one still open, H29485893, moves
# from row 8 up to row 1), row 3 becomes a new "corregir" entry (with
# blank follow-up cells C3:E3), and the three old "corregir" rows
# (11-13) are cleared out since they've been resolved/removed.

$ws.Range("A1").Value = "H29485893"
$ws.Range("B1").Value = "terminado"

$ws.Range("A2").Value = "H29498961"
$ws.Range("B2").Value = "terminado"

$ws.Range("A3").Value = "H29531290"
$ws.Range("B3").Value = "corregir"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""

$ws.Range("A4").Value = "H29552494"
$ws.Range("B4").Value = "terminado"

$ws.Range("A5").Value = "H29552791"
$ws.Range("B5").Value = "terminado"

$ws.Range("A6").Value = "H29552908"
$ws.Range("B6").Value = "terminado"

$ws.Range("A7").Value = "H29568474"
$ws.Range("B7").Value = "terminado"

$ws.Range("A8").Value = "H29617552"
$ws.Range("B8").Value = "terminado"

$ws.Range("A9").Value = "H29634433"
$ws.Range("B9").Value = "terminado"

$ws.Range("A10").Value = "H29674926"
$ws.Range("B10").Value = "terminado"

# Rows 11-13 no longer have pending "corregir" tickets - clear them.
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = ""

$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""

$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
